$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feature names (column A) and importance values (column B)
# Rows 2-10 reordered/updated per new feature importance ranking

$ws.Range("A2").Value = "MACD"
$ws.Range("B2").Value = 0.3930817313426748

$ws.Range("A3").Value = "RSI"
$ws.Range("B3").Value = 0.3080310985219448

$ws.Range("A4").Value = "Signal_line"
$ws.Range("B4").Value = 0.1320233135593225

$ws.Range("A5").Value = "VIX_short"
$ws.Range("B5").Value = 0.04050396027363317

$ws.Range("A6").Value = "close_short"
$ws.Range("B6").Value = 0.03403665931781309

$ws.Range("A7").Value = "close_long"
$ws.Range("B7").Value = 0.03071072798458908

$ws.Range("A8").Value = "VIX"
$ws.Range("B8").Value = 0.02369486540704165

$ws.Range("A9").Value = "VIX_long"
$ws.Range("B9").Value = 0.02211115837664275

$ws.Range("A10").Value = "fedrate"
$ws.Range("B10").Value = 0.01580648521633805
